# Remove the trailing "Ver no Jupiter ..." and "(c) 2020 ..." site-chrome
# paragraphs (plus the blank paragraph separating them from the
# bibliography), left over from a Jekyll site rebuild. The bibliography's
# last entry ("Janeiro: Editora Interciencia , 2004.") and the final blank
# paragraph before the page break are left untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by content rather than a
# hard-coded index, so the script is resilient to minor shifts.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Ver no Jupiter*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 1) {
    # The blank paragraph immediately preceding "Ver no Jupiter ..."
    $pBlank = $d.Paragraphs($targetIndex - 1)
    # The "(c) 2020 ... Creative Commons Attribution" paragraph right after it
    $pCopyright = $d.Paragraphs($targetIndex + 1)

    $start = $pBlank.Range.Start
    $end = $pCopyright.Range.End

    $d.Range($start, $end).Delete()
}
